# Updates the crypto price/volume columns (D, E) to the latest snapshot values.
# A handful of D-column cells are numeric-looking text (e.g. "5.360", "0.9730",
# "0.000006806") whose trailing zeros / scale Excel would mangle if it auto-
# converted them to real numbers on assignment, so those are written with a
# leading apostrophe to force them to stay literal text (matches how the
# original workbook stores every value in these columns: as text).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.356.51'
$ws.Range("D3").Value = '1.719.48'
$ws.Range("D4").Value = '0.9987'
$ws.Range("D5").Value = '239.42'
$ws.Range("D6").Value = '0.9997'
$ws.Range("D7").Value = '0.4719'
$ws.Range("D9").Value = '0.06216'
$ws.Range("D10").Value = '1.712.20'
$ws.Range("D11").Value = '0.07073'
$ws.Range("D12").Value = '15.26'
$ws.Range("D13").Value = '0.5924'
$ws.Range("D14").Value = '4.411'
$ws.Range("D15").Value = '76.48'
$ws.Range("D17").Value = '0.9995'
$ws.Range("D18").Value = '26.350.90'
$ws.Range("D19").Value = "'0.000006806"
$ws.Range("D21").Value = '1.930.88'
$ws.Range("D22").Value = '4.549'
$ws.Range("D23").Value = '8.808'
$ws.Range("D24").Value = "'5.360"
$ws.Range("D25").Value = '135.42'
$ws.Range("D26").Value = "'15.20"
$ws.Range("D27").Value = '1.404'
$ws.Range("D28").Value = '1.765'
$ws.Range("D29").Value = '106.96'
$ws.Range("D30").Value = '4.051'
$ws.Range("D31").Value = '3.689'
$ws.Range("D32").Value = '0.07716'
$ws.Range("D33").Value = '0.04432'
$ws.Range("D34").Value = '2.611'
$ws.Range("D35").Value = '0.6224'
$ws.Range("D36").Value = "'0.9730"
$ws.Range("D37").Value = "'0.9360"
$ws.Range("D38").Value = '115.25'
$ws.Range("D39").Value = '2.412'
$ws.Range("D41").Value = '1.907'
$ws.Range("D42").Value = '0.01468'
$ws.Range("D43").Value = '5.298'
$ws.Range("D44").Value = '0.3808'
$ws.Range("D45").Value = '0.1152'
$ws.Range("D46").Value = "'6.260"
$ws.Range("D47").Value = '0.05288'
$ws.Range("D48").Value = '30.52'
$ws.Range("D49").Value = '7.633'
$ws.Range("D51").Value = '0.3383'
$ws.Range("E2").Value = '  +3.86%  '
$ws.Range("E3").Value = '  +3.33%  '
$ws.Range("E5").Value = '  +1.46%  '
$ws.Range("E7").Value = '  -1.54%  '
$ws.Range("E8").Value = '  +1.04%  '
$ws.Range("E9").Value = '  +1.11%  '
$ws.Range("E10").Value = '  +2.96%  '
$ws.Range("E11").Value = '  -0.08%  '
$ws.Range("E12").Value = '  +3.90%  '
$ws.Range("E13").Value = '  +0.19%  '
$ws.Range("E14").Value = '  +0.71%  '
$ws.Range("E15").Value = '  +3.01%  '
$ws.Range("E16").Value = '  -0.02%  '
$ws.Range("E17").Value = '  -0.04%  '
$ws.Range("E18").Value = '  +3.84%  '
$ws.Range("E19").Value = '  +1.22%  '
$ws.Range("E20").Value = '  +1.20%  '
$ws.Range("E21").Value = '  +3.14%  '
$ws.Range("E22").Value = '  +2.52%  '
$ws.Range("E23").Value = '  +1.81%  '
$ws.Range("E24").Value = '  +0.63%  '
$ws.Range("E25").Value = '  +1.44%  '
$ws.Range("E26").Value = '  +0.83%  '
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("E28").Value = '  +4.35%  '
$ws.Range("E29").Value = '  +2.47%  '
$ws.Range("E30").Value = '  +1.97%  '
$ws.Range("E31").Value = '  +1.93%  '
$ws.Range("E32").Value = '  +0.87%  '
$ws.Range("E33").Value = '  +0.93%  '
$ws.Range("E34").Value = '  +0.38%  '
$ws.Range("E35").Value = '  +2.34%  '
$ws.Range("E36").Value = '  +3.28%  '
$ws.Range("E37").Value = '  +9.82%  '
$ws.Range("E38").Value = '  +17.01%  '
$ws.Range("E39").Value = '  -8.12%  '
$ws.Range("E40").Value = '  -0.03%  '
$ws.Range("E41").Value = '  +4.82%  '
$ws.Range("E42").Value = '  -2.18%  '
$ws.Range("E43").Value = '  +13.02%  '
$ws.Range("E44").Value = '  +1.25%  '
$ws.Range("E45").Value = '  +3.13%  '
$ws.Range("E46").Value = '  +0.81%  '
$ws.Range("E47").Value = '  +0.70%  '
$ws.Range("E48").Value = '  +3.37%  '
$ws.Range("E49").Value = '  +4.66%  '
$ws.Range("E50").Value = '  +0.86%  '
$ws.Range("E51").Value = '  +1.16%  '
